$wb = $excel.ActiveWorkbook

# --- feb2025 sheet: mark "pago2" (column D) as paid for several rows ---
$wsFeb = $wb.Worksheets.Item("feb2025")
$wsFeb.Range("D3").Value = 65000
$wsFeb.Range("D4").Value = 65000
$wsFeb.Range("D5").Value = 65000
$wsFeb.Range("D6").Value = 65000
$wsFeb.Range("D11").Value = 65000
$wsFeb.Range("D16").Value = 65000
$wsFeb.Range("D17").Value = 65000
$wsFeb.Range("D23").Value = 65000
$wsFeb.Range("D24").Value = 65000
$wsFeb.Activate()
$wsFeb.Range("D7").Select()

# --- mar2025 sheet: mark "pago1" (column C) as paid for row 16 ---
$wsMar = $wb.Worksheets.Item("mar2025")
$wsMar.Range("C16").Value = 65000
$wsMar.Activate()
$wsMar.Range("C17").Select()

# restore original active sheet (feb2025, which was tabSelected="true")
$wsFeb.Activate()
